$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column A: was stored as text "79174445" -> convert to numeric 79174445
$ws.Range("A7").Value = 79174445

# Append new row 8 with the redemption: phone 71277628 redeems 80 points
# Phone numbers are kept as text (like the other phone cells), so force
# text formatting before assignment, then drop back to the Normal style
# so no extra formatting is visibly applied to the cell.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "71277628"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = 80

$ws.Range("C8").Value = "2025-08-18T16:53:26"
